$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("47:48").Insert()

$data = New-Object 'object[,]' 1,20
$vals = @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 45028, 13, "Fruta", 100104, "Frutos de pepita", 100104003, "Membrillo", "Champion", "Primera", 14, 230000, 240000, 235714, "`$/bins (450 kilos)", "Región de O'Higgins", 524, 450)
for ($i = 0; $i -lt 20; $i++) {
    $data[0, $i] = $vals[$i]
}
$ws.Range("A47:T47").Value = $data
$v = $ws.Range("A47").Value()
Write-Host ("A47 after 2D array set = " + [string]$v)
$v2 = $ws.Range("L47").Value()
Write-Host ("L47 after 2D array set = " + [string]$v2)
$v3 = $ws.Range("T47").Value()
Write-Host ("T47 after 2D array set = " + [string]$v3)
